# Replace every occurrence of the acronym "OIE" with "WOAH" across all
# worksheets (the World Organisation for Animal Health rebranded from
# "OIE" to "WOAH"). This covers the narrative text cells (Sheet 1) and the
# reference-citation cells (References), while leaving lower-case
# occurrences inside URLs (e.g. "https://www.oie.int/...") untouched,
# exactly as in the source diff.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    # Replace("OIE" -> "WOAH"), substring match (xlPart = 2), by rows
    # (xlByRows = 1), case-sensitive match so "oie.int" URLs are left alone.
    [void]$usedRange.Replace("OIE", "WOAH", 2, 1, $true)
}
